$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.573.88'
$ws.Range("E2").Value = '  -0.23%  '
$ws.Range("D3").Value = '3.334.51'
$ws.Range("E3").Value = '  +0.48%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '580.88'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.75%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '175.89'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.20%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("E8").Value = '  +0.10%  '
$ws.Range("D9").Value = '3.334.37'
$ws.Range("E9").Value = '  +0.74%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.179'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.33%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.577'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.12%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '45.42'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.03%  '
$ws.Range("E13").Value = '  -2.32%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '668.03'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +4.27%  '
$ws.Range("D15").Value = '3.879.27'
$ws.Range("E15").Value = '  +0.76%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '8.42'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.27%  '
$ws.Range("D17").Value = '67.771.40'
$ws.Range("E17").Value = '  -0.12%  '
$ws.Range("E18").Value = '  -0.46%  '
$ws.Range("D19").Value = '3.339.28'
$ws.Range("E19").Value = '  +0.47%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.40'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.51%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.96'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.21%  '
$ws.Range("E22").Value = '  -0.81%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.49'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +9.51%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '17.06'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.59%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '99.42'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.49%  '
$ws.Range("E26").Value = '  -3.26%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.67'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -4.39%  '
$ws.Range("E28").Value = '  -2.98%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '33.64'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.66%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.39'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +11.01%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.45'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.37%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '576.58'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -5.00%  '
$ws.Range("E33").Value = '  +0.31%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.104'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.20%  '
$ws.Range("E35").Value = '  +0.09%  '
$ws.Range("D36").Value = '3.703.58'
$ws.Range("E36").Value = '  -5.69%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '56.66'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.62%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.36'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -6.40%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '34.53'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.93%  '
$ws.Range("E40").Value = '  +2.16%  '
$ws.Range("E41").Value = '  -2.69%  '
$ws.Range("E42").Value = '  -4.87%  '
$ws.Range("D43").Value = '0.0₃0670'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.334'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.89%  '
$ws.Range("E45").Value = '  -1.64%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0406'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.15%  '
$ws.Range("E47").Value = '  +2.01%  '
$ws.Range("E48").Value = '  -0.39%  '
$ws.Range("E49").Value = '  -0.12%  '
$ws.Range("E50").Value = '  +1.13%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '128.31'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.24%  '
